$d = $word.ActiveDocument

# Update the bookmark-style placeholder text in the first paragraph.
$range = $d.Paragraphs.Item(1).Range
$found = $range.Find.Execute(
    "**ID__AFFARS_mp_5306_502_topic_3__ID**", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "**ID__AFFARS_MP_5306_502_2__ID**", 2)

# Remove the now-orphaned trailing space run that used to follow the
# placeholder text (the Find/Replace above merges runs, leaving a single
# run whose text still ends with the old separating space).
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$trailRange = $d.Range($r.End - 2, $r.End - 1)
if ($trailRange.Text -eq " ") {
    $trailRange.Delete()
}

# Give the paragraph a (currently invisible/unset) border, matching the
# border used by sibling paragraphs in this document, and widen the left
# indent to match.
$p1 = $d.Paragraphs.Item(1)
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25
